$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: task "Integration for the whole project" (Third phase coding and testing) ---
# To-Do deadline (B16), Done date for khaled (D16) and Dapsara (F16) move from 4/7/2023 to 6/4/2023.
$ws.Range("B16").Value = 45081
$ws.Range("D16").Value = 45081
$ws.Range("F16").Value = 45081

# C16 used to hold the text "on going " -- it now becomes a real "Done" date (6/4/2023),
# so give it the same number-format/style as the other date cells in the row first.
$ws.Range("B16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C16").Value = 45081

# --- Row 17: new task "Forth Phase for softwae and coding" ---
$ws.Range("A17").Value = "Forth Phase for softwae and coding"

$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B17").Value = 45091
$ws.Range("C17").Value = 45091

# --- Cosmetic: move the active selection to B14 ---
$ws.Range("B14").Select() | Out-Null
